# Scheduled runner update: refresh Leve profit calculations (currentAveragePrice*,
# LevePrice*, LeveProfit*) across ALC/ARM/BSM/CUL/GSM sheets with latest market data.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 137
$ws.Range("H137").Value = 1350.625
$ws.Range("I137").Value = 1222.6666
$ws.Range("J137").Value = 1515.1428
$ws.Range("K137").Value = 3667.9998
$ws.Range("L137").Value = 4545.428400000001
$ws.Range("M137").Value = -1117.9998
$ws.Range("N137").Value = -9645.428400000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 2227.6924
$ws.Range("I2").Value = 2398
$ws.Range("J2").Value = 1955.2
$ws.Range("K2").Value = 2398
$ws.Range("L2").Value = 1955.2
$ws.Range("M2").Value = -2285
$ws.Range("N2").Value = -2181.2
# row 74
$ws.Range("H74").Value = 914.1489
$ws.Range("I74").Value = 641.2059
$ws.Range("J74").Value = 1628
$ws.Range("K74").Value = 641.2059
$ws.Range("L74").Value = 1628
$ws.Range("M74").Value = 232.7941
$ws.Range("N74").Value = -3376
# row 77
$ws.Range("H77").Value = 914.1489
$ws.Range("I77").Value = 641.2059
$ws.Range("J77").Value = 1628
$ws.Range("K77").Value = 3206.0295
$ws.Range("L77").Value = 8140
$ws.Range("M77").Value = 1161.9705
$ws.Range("N77").Value = -16876
# row 116
$ws.Range("H116").Value = 2227.6924
$ws.Range("I116").Value = 2398
$ws.Range("J116").Value = 1955.2
$ws.Range("K116").Value = 2398
$ws.Range("L116").Value = 1955.2
$ws.Range("M116").Value = -104
$ws.Range("N116").Value = -6543.2

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 2227.6924
$ws.Range("I3").Value = 2398
$ws.Range("J3").Value = 1955.2
$ws.Range("K3").Value = 2398
$ws.Range("L3").Value = 1955.2
$ws.Range("M3").Value = -2284
$ws.Range("N3").Value = -2183.2
# row 86
$ws.Range("H86").Value = 1383.1305
$ws.Range("I86").Value = 1381.5238
$ws.Range("J86").Value = 1400
$ws.Range("K86").Value = 1381.5238
$ws.Range("L86").Value = 1400
$ws.Range("M86").Value = -258.5237999999999
$ws.Range("N86").Value = -3646
# row 89
$ws.Range("H89").Value = 1383.1305
$ws.Range("I89").Value = 1381.5238
$ws.Range("J89").Value = 1400
$ws.Range("K89").Value = 6907.619
$ws.Range("L89").Value = 7000
$ws.Range("M89").Value = -1291.619
$ws.Range("N89").Value = -18232
# row 132
$ws.Range("H132").Value = 40125
$ws.Range("J132").Value = 40125
$ws.Range("L132").Value = 40125
$ws.Range("N132").Value = -50245

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 68
$ws.Range("H68").Value = 1769.8607
$ws.Range("I68").Value = 1382.325
$ws.Range("J68").Value = 2167.3333
$ws.Range("K68").Value = 4146.975
$ws.Range("L68").Value = 6501.999899999999
$ws.Range("M68").Value = -3335.975
$ws.Range("N68").Value = -8123.999899999999
# row 71
$ws.Range("H71").Value = 1769.8607
$ws.Range("I71").Value = 1382.325
$ws.Range("J71").Value = 2167.3333
$ws.Range("K71").Value = 12440.925
$ws.Range("L71").Value = 19505.9997
$ws.Range("M71").Value = -8384.925000000001
$ws.Range("N71").Value = -27617.9997
# row 93
$ws.Range("H93").Value = 7900
$ws.Range("J93").Value = 7900
$ws.Range("L93").Value = 23700
$ws.Range("N93").Value = -27444
# row 95
$ws.Range("H95").Value = 6565.4
$ws.Range("J95").Value = 6565.4
$ws.Range("L95").Value = 19696.2
$ws.Range("N95").Value = -23814.2
# row 98
$ws.Range("H98").Value = 212.8
$ws.Range("I98").Value = 183.33333
$ws.Range("J98").Value = 257
$ws.Range("K98").Value = 549.99999
$ws.Range("L98").Value = 771
$ws.Range("M98").Value = 948.00001
$ws.Range("N98").Value = -3767
# row 111
$ws.Range("H111").Value = 3084.111
$ws.Range("I111").Value = 3206.75
$ws.Range("J111").Value = 2986
$ws.Range("K111").Value = 9620.25
$ws.Range("L111").Value = 8958
$ws.Range("M111").Value = -6553.25
$ws.Range("N111").Value = -15092
# row 112
$ws.Range("H112").Value = 4964.615
$ws.Range("I112").Value = 2750
$ws.Range("J112").Value = 5367.273
$ws.Range("K112").Value = 8250
$ws.Range("L112").Value = 16101.819
$ws.Range("M112").Value = -7142
$ws.Range("N112").Value = -18317.819
# row 120
$ws.Range("H120").Value = 10500
$ws.Range("I120").Value = 7333.3335
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 22000.0005
$ws.Range("L120").Value = 60000
$ws.Range("M120").Value = -17162.0005
$ws.Range("N120").Value = -69676
# row 121
$ws.Range("H121").Value = 14620750
$ws.Range("I121").Value = 650
$ws.Range("J121").Value = 15432978
$ws.Range("K121").Value = 1950
$ws.Range("L121").Value = 46298934
$ws.Range("M121").Value = -640
$ws.Range("N121").Value = -46301554
# row 122
$ws.Range("H122").Value = 2277900.8
$ws.Range("I122").Value = 509.5
$ws.Range("J122").Value = 5010770.5
$ws.Range("K122").Value = 4585.5
$ws.Range("L122").Value = 45096934.5
$ws.Range("M122").Value = -2135.5
$ws.Range("N122").Value = -45101834.5
# row 123
$ws.Range("H123").Value = 2033
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 2033
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 6099
$ws.Range("N123").Value = -10999
# row 124
$ws.Range("H124").Value = 1864.5
$ws.Range("I124").Value = 890
$ws.Range("J124").Value = 2514.1667
$ws.Range("K124").Value = 2670
$ws.Range("L124").Value = 7542.500100000001
$ws.Range("M124").Value = 2240
$ws.Range("N124").Value = -17362.5001
# row 125
$ws.Range("H125").Value = 9371.5
$ws.Range("I125").Value = 1500
$ws.Range("J125").Value = 11995.333
$ws.Range("K125").Value = 4500
$ws.Range("L125").Value = 35985.999
$ws.Range("M125").Value = 420
$ws.Range("N125").Value = -45825.999
# row 126
$ws.Range("H126").Value = 4655.3335
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4655.3335
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 13966.0005
$ws.Range("N126").Value = -23846.0005
# row 127
$ws.Range("H127").Value = 2033
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 2033
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 6099
$ws.Range("N127").Value = -16019
# row 128
$ws.Range("H128").Value = 150000
$ws.Range("I128").Value = 150000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 450000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -445020
# row 129
$ws.Range("H129").Value = 1255.0714
$ws.Range("I129").Value = 1056.3334
$ws.Range("J129").Value = 1404.125
$ws.Range("K129").Value = 3169.0002
$ws.Range("L129").Value = 4212.375
$ws.Range("M129").Value = 1830.9998
$ws.Range("N129").Value = -14212.375
# row 130
$ws.Range("H130").Value = 4155.75
$ws.Range("I130").Value = 3030
$ws.Range("J130").Value = 4531
$ws.Range("K130").Value = 9090
$ws.Range("L130").Value = 13593
$ws.Range("M130").Value = -4070
$ws.Range("N130").Value = -23633
# row 131
$ws.Range("H131").Value = 932.6875
$ws.Range("I131").Value = 397.66666
$ws.Range("J131").Value = 1056.1538
$ws.Range("K131").Value = 1192.99998
$ws.Range("L131").Value = 3168.4614
$ws.Range("M131").Value = 3847.00002
$ws.Range("N131").Value = -13248.4614
# row 132
$ws.Range("H132").Value = 579.3333
$ws.Range("I132").Value = 557.5714
$ws.Range("J132").Value = 590.2143
$ws.Range("K132").Value = 5018.1426
$ws.Range("L132").Value = 5311.928699999999
$ws.Range("M132").Value = -2488.1426
$ws.Range("N132").Value = -10371.9287
# row 133
$ws.Range("H133").Value = 4862.222
$ws.Range("I133").Value = 4465.7144
$ws.Range("J133").Value = 6250
$ws.Range("K133").Value = 13397.1432
$ws.Range("L133").Value = 18750
$ws.Range("M133").Value = -8337.143199999999
$ws.Range("N133").Value = -28870
# row 134
$ws.Range("H134").Value = 2623.4
$ws.Range("I134").Value = 1669.7059
$ws.Range("J134").Value = 4650
$ws.Range("K134").Value = 5009.1177
$ws.Range("L134").Value = 13950
$ws.Range("M134").Value = 60.88230000000021
$ws.Range("N134").Value = -24090
# row 136
$ws.Range("H136").Value = 14709795
$ws.Range("I136").Value = 25001290
$ws.Range("J136").Value = 7657.143
$ws.Range("K136").Value = 75003870
$ws.Range("L136").Value = 22971.429
$ws.Range("M136").Value = -74998770
$ws.Range("N136").Value = -33171.429
# row 137
$ws.Range("H137").Value = 20881034
$ws.Range("I137").Value = 38470180
$ws.Range("J137").Value = 93860.55
$ws.Range("K137").Value = 115410540
$ws.Range("L137").Value = 281581.65
$ws.Range("M137").Value = -115405440
$ws.Range("N137").Value = -291781.65
# row 138
$ws.Range("H138").Value = 2714
$ws.Range("I138").Value = 1520
$ws.Range("J138").Value = 5997.5
$ws.Range("K138").Value = 4560
$ws.Range("L138").Value = 17992.5
$ws.Range("M138").Value = 580
$ws.Range("N138").Value = -28272.5
# row 139
$ws.Range("H139").Value = 1581.6111
$ws.Range("I139").Value = 1498.1765
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 4494.529500000001
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 645.4704999999994
$ws.Range("N139").Value = -19280
# row 140
$ws.Range("H140").Value = 1832.4
$ws.Range("I140").Value = 1134.4615
$ws.Range("J140").Value = 3128.5715
$ws.Range("K140").Value = 3403.3845
$ws.Range("L140").Value = 9385.7145
$ws.Range("M140").Value = 1776.6155
$ws.Range("N140").Value = -19745.7145
# row 141
$ws.Range("H141").Value = 4122.357
$ws.Range("I141").Value = 3348.3333
$ws.Range("J141").Value = 8766.5
$ws.Range("K141").Value = 10044.9999
$ws.Range("L141").Value = 26299.5
$ws.Range("M141").Value = -4864.999899999999
$ws.Range("N141").Value = -36659.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 2761.9614
$ws.Range("I102").Value = 2534.2
$ws.Range("J102").Value = 3072.5454
$ws.Range("K102").Value = 2534.2
$ws.Range("L102").Value = 3072.5454
$ws.Range("M102").Value = -912.1999999999998
$ws.Range("N102").Value = -6316.5454
